$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-19 Friday" "2025-12-20 Saturday"

Replace-Text "662×2=1324" "793×6=4758"
Replace-Text "251×6=1506" "745×6=4470"
Replace-Text "409×6=2454" "265×8=2120"
Replace-Text "477×2=954" "745×6=4470"
Replace-Text "156×7=1092" "261×9=2349"

Replace-Text "787×9=7083" "790×4=3160"
Replace-Text "301×4=1204" "804×8=6432"
Replace-Text "140×4=560" "519×4=2076"
Replace-Text "468×5=2340" "380×3=1140"
Replace-Text "286×4=1144" "447×5=2235"

Replace-Text "151×7=1057" "321×7=2247"
Replace-Text "186×4=744" "772×7=5404"
Replace-Text "790×3=2370" "163×4=652"
Replace-Text "596×3=1788" "436×8=3488"
Replace-Text "737×3=2211" "914×7=6398"

Replace-Text "385×9=3465" "704×4=2816"
Replace-Text "589×3=1767" "351×4=1404"
Replace-Text "471×7=3297" "840×3=2520"
Replace-Text "679×7=4753" "840×3=2520"
Replace-Text "478×7=3346" "408×4=1632"

Replace-Text "957×8=7656" "920×7=6440"
Replace-Text "169×4=676" "392×5=1960"
Replace-Text "758×3=2274" "194×3=582"
Replace-Text "699×9=6291" "323×5=1615"
Replace-Text "598×3=1794" "470×8=3760"
